# Fix mojibake "(ÂºC)" -> "(ºC)" in the temperature-range table header cells.
$d = $word.ActiveDocument

$find = [string][char]0x00C2 + [string][char]0x00BA + "C"
$replace = [string][char]0x00BA + "C"

$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
